$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(196, 1).Value = 45998
$ws.Cells.Item(196, 2).Value = "四方坪站充电量(kw)"
$ws.Cells.Item(196, 3).Value = 608.52500000000009
$ws.Cells.Item(196, 4).Value = 1021.9850000000001
$ws.Cells.Item(196, 5).Value = 453.50200000000001
$ws.Cells.Item(196, 6).Value = 485.76
$ws.Cells.Item(196, 7).Value = 85.649000000000001
$ws.Cells.Item(196, 8).Value = 552.61099999999999
$ws.Cells.Item(196, 9).Value = 492.33900000000006
$ws.Cells.Item(196, 10).Value = 170.24699999999999
$ws.Cells.Item(196, 11).Value = 113.114
$ws.Cells.Item(196, 12).Value = 115.52800000000002
$ws.Cells.Item(196, 13).Value = 216.01700000000002
$ws.Cells.Item(196, 14).Value = 212.95099999999999
$ws.Cells.Item(196, 15).Value = 498.18999999999994
$ws.Cells.Item(196, 16).Value = 1138.4869999999996
$ws.Cells.Item(196, 17).Value = 360.78999999999996
$ws.Cells.Item(196, 18).Value = 449.399
$ws.Cells.Item(196, 19).Value = 333.6
$ws.Cells.Item(196, 20).Value = 409.45300000000003
$ws.Cells.Item(196, 21).Value = 148.14400000000003
$ws.Cells.Item(196, 22).Value = 40.54
$ws.Cells.Item(196, 23).Value = 76.16
$ws.Cells.Item(196, 24).Value = 76.599999999999994
$ws.Cells.Item(196, 25).Value = 95.62
$ws.Cells.Item(196, 26).Value = 84.960000000000008

$ws.Cells.Item(197, 1).Value = 45998
$ws.Cells.Item(197, 2).Value = "高岭站充电量(kw)"
$ws.Cells.Item(197, 3).Value = 293.08999999999997
$ws.Cells.Item(197, 4).Value = 597.178
$ws.Cells.Item(197, 5).Value = 98.460999999999984
$ws.Cells.Item(197, 6).Value = 60.938000000000002
$ws.Cells.Item(197, 7).Value = 0
$ws.Cells.Item(197, 8).Value = 75.085999999999999
$ws.Cells.Item(197, 9).Value = 387.47599999999994
$ws.Cells.Item(197, 10).Value = 142.14599999999999
$ws.Cells.Item(197, 11).Value = 233.196
$ws.Cells.Item(197, 12).Value = 231.05900000000003
$ws.Cells.Item(197, 13).Value = 142.83399999999997
$ws.Cells.Item(197, 14).Value = 302.72699999999998
$ws.Cells.Item(197, 15).Value = 447.58100000000002
$ws.Cells.Item(197, 16).Value = 561.63099999999986
$ws.Cells.Item(197, 17).Value = 302.86700000000002
$ws.Cells.Item(197, 18).Value = 240.953
$ws.Cells.Item(197, 19).Value = 177.85299999999998
$ws.Cells.Item(197, 20).Value = 93.304000000000002
$ws.Cells.Item(197, 21).Value = 0
$ws.Cells.Item(197, 22).Value = 66.164000000000001
$ws.Cells.Item(197, 23).Value = 27.688000000000002
$ws.Cells.Item(197, 24).Value = 88.935000000000002
$ws.Cells.Item(197, 25).Value = 50.058
$ws.Cells.Item(197, 26).Value = 32.773000000000003

[void]$ws.Range("G203").Select()
